$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Oeffnungszeiten" cells carried a stray legacy font variant;
# normalize them onto the sheet's regular body font (matches the rest of
# column D/F) now that new rows are being added below them.
foreach ($addr in @("D24", "D26", "D28", "F29", "D32")) {
    $ws.Range($addr).Font.Name = "Arial"
    $ws.Range($addr).Font.Size = 12
}

# E31 ("Stosszeit" for Bio-Laden Rosmarin) keeps its time-of-day format.
$ws.Range("E31").NumberFormat = "hh:mm:ss\ AM/PM"

# Row 34: "Stoff-Art" fabric shop
$ws.Range("B34").Value = "Stoff-Art"
$ws.Range("C34").Value = "Eisfeld 3, 99423 Weimar"
$ws.Range("D34").Value = "Mo – Fr 10:00 – 18:00, Sa 10:00 – 14:00"
$ws.Range("E34").Value = 0.5
$ws.Range("E34").NumberFormat = "hh:mm:ss\ AM/PM"
$ws.Range("F34").Value = "4.8/5 (46)"
$ws.Range("G34").Value = "Stoffgeschäft"

# Row 35: "Verfilzt & Zugenäht" fabric shop
$ws.Range("B35").Value = "Verfilzt & Zugenäht"
$ws.Range("C35").Value = "Jakobstraße 2, 99423 Weimar"
$ws.Range("D35").Value = "Mo – Fr 11:00 – 17:00, Sa 11:00 – 15:00"
$ws.Range("E35").Value = "n.a."
$ws.Range("F35").Value = "4.9/5 (39)"
$ws.Range("G35").Value = "Stoffgeschäft"

# Selection matches the author's final saved view state.
$ws.Range("E19").Select()
